$d = $word.ActiveDocument

# 1) "Fixing sign location": the signature block used a hard-coded city
#    name ("Jakarta") in front of the date placeholder. Replace the
#    literal word "Jakarta" with the "${authority_location}" merge-field
#    placeholder everywhere it appears (the trailing punctuation / spacing
#    that follows it, e.g. ",  " or ", ", is left untouched).
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Jakarta", $false, $true, $false, $false, $false, $true, 1, $false, '${authority_location}', 2)

# 2) Remove three of the superfluous empty paragraphs that were sitting
#    right after the last table (they were accidental duplicates left
#    over in the template).
$tbl = $d.Tables.Item($d.Tables.Count)
$extra = $d.Range($tbl.Range.End, $tbl.Range.End + 3)
$extra.Delete()
